# Updated cryptos list data (price and 1h volume change columns),
# plus a rank swap between Litecoin and Polygon (rows 15/16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.893.10'

# Row 3
$ws.Range('D3').Value = '1.665.63'
$ws.Range('E3').Value = '  +0.79%  '

# Row 4
$ws.Range('E4').Value = '  -0.20%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.532'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.93%  '

# Row 7
$ws.Range('E7').Value = '  -0.17%  '

# Row 8
$ws.Range('E8').Value = '  +0.66%  '

# Row 9
$ws.Range('E9').Value = '  +1.05%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.36%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0897'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.82%  '

# Row 12
$ws.Range('D12').Value = '1.899.34'
$ws.Range('E12').Value = '  +0.66%  '

# Row 13
$ws.Range('D13').Value = '1.668.73'
$ws.Range('E13').Value = '  +0.51%  '

# Row 14
$ws.Range('E14').Value = '  +0.02%  '

# Row 15
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '66.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.34%  '

# Row 16
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.525'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.19%  '

# Row 17
$ws.Range('D17').Value = '26.888.88'
$ws.Range('E17').Value = '  -0.52%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.66%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.14%  '

# Row 20
$ws.Range('E20').Value = '  +0.33%  '

# Row 21
$ws.Range('E21').Value = '  -0.26%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.77%  '

# Row 23
$ws.Range('E23').Value = '  -1.03%  '

# Row 24
$ws.Range('E24').Value = '  -1.24%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
$ws.Range('E26').Value = '  -0.15%  '

# Row 27
$ws.Range('E27').Value = '  +1.63%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.77%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '

# Row 30
$ws.Range('E30').Value = '  +0.09%  '

# Row 31
$ws.Range('E31').Value = '  +0.25%  '

# Row 32
$ws.Range('E32').Value = '  +2.03%  '

# Row 33
$ws.Range('D33').Value = '1.458.02'
$ws.Range('E33').Value = '  -4.66%  '

# Row 34
$ws.Range('E34').Value = '  +3.79%  '

# Row 35
$ws.Range('E35').Value = '  +3.25%  '

# Row 36
$ws.Range('E36').Value = '  -0.33%  '

# Row 37
$ws.Range('E37').Value = '  +0.50%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.901'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.94%  '

# Row 39
$ws.Range('E39').Value = '  +0.40%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.31%  '

# Row 41
$ws.Range('E41').Value = '  -0.25%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.41%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.977'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.03%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.46%  '

# Row 45
$ws.Range('D45').Value = '1.808.59'
$ws.Range('E45').Value = '  +0.79%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.781'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.80%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.32%  '

# Row 48
$ws.Range('E48').Value = '  +0.89%  '

# Row 49
$ws.Range('E49').Value = '  -1.60%  '

# Row 50
$ws.Range('E50').Value = '  +4.60%  '

# Row 51
$ws.Range('E51').Value = '  +0.52%  '

